# The document has a run of "pseudocode" paragraphs describing an
# application's logic. Previously, the paragraph right after the
# "...call api call function save to variable" line was a bare empty
# paragraph (<w:p/>). It needs to become a complete line of pseudocode,
# and a fresh empty (tab-only) paragraph needs to follow it, matching the
# pattern used by the rest of the document (each paragraph starts with a
# tab to show nesting under its "Declare function" parent).

$d = $word.ActiveDocument

# Find the anchor line so we don't depend on hard-coded paragraph
# indices, then grab the paragraph immediately following it -- that is
# the empty paragraph we need to fill in.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("call function save to variable", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchorRange.Paragraphs(1)
$targetPara = $anchorPara.Next()

# Fill the previously-empty paragraph with the new pseudocode line
# (indented with a leading tab, like its sibling lines).
$targetRange = $targetPara.Range
$targetRange.InsertAfter([char]9 + "Create variable for status of api call from json data")

# Insert a brand-new empty paragraph (containing just a tab, like the
# other spacer paragraphs in the document) right after it.
$targetRange.InsertParagraphAfter()
$spacerPara = $targetPara.Next()
$spacerPara.Range.InsertAfter([char]9)
